# Update existing row 2 (Todoist -> Synchroteam) and add new row 3 (Yeastar)
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("API_Results")

# Row 2: replace Todoist entry with Synchroteam entry
$ws.Range("A2").Value = "Synchroteam"
$ws.Range("B2").Value = "Field Service Management"
$ws.Range("C2").Value = "Scheduling & Dispatching"
$ws.Range("D2").Value = "Yes (14 days - Full) - 95%"
$ws.Range("E2").Value = "REST API"
$ws.Range("F2").Value = "Yes - /User/Send, /User/List"
$ws.Range("G2").Value = "https://www.synchroteam.com"
$ws.Range("H2").Value = "No credit card required"

# Row 3: new Yeastar entry
$ws.Range("A3").Value = "Yeastar"
$ws.Range("B3").Value = "Unified Communications"
$ws.Range("C3").Value = "Business Phone Systems (PBX)"
$ws.Range("D3").Value = "Yes (30 days - Full) - 90%"
$ws.Range("E3").Value = "REST API"
$ws.Range("F3").Value = "Yes - /extensionlist/query, /extension/update"
$ws.Range("G3").Value = "https://www.yeastar.com"
$ws.Range("H3").Value = "No credit card required"
